$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.67%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'36.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.87%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.035"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.12%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07914"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.43%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.126"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-2.63%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.966"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.83%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9230"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.53%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.09811"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.99%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1863"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.20%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08972"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.79%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03598"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.28%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09925"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.001436"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-2.75%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005607"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.50%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.480"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.53%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.139"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.95%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'8.97%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3424"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.31%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1335"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.61%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.174"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'4.79%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2244"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.89%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04563"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.14%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-1.32%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004830"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-7.30%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001297"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.49%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0004741"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'74.28%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01864"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.81%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04904"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.94%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007783"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.68%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1401"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.83%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007713"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.79%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002211"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'3.23%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01125"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'11.34%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006413"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.14%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.34%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'0.13%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'51.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'42.98%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.001897"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-29.54%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.00002095"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.34%"
$ws.Range("E51").Style = "Normal"
